$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh NN hyperparameter-tuning metrics for the new-laptop training run ---
# Rows 2-32 hold per-NN-config metrics (cvMSE/cvRMSE/cvMAE stay 0; Train R^2,
# Train Adjusted R^2, Train MSE, Train RMSE, Train MAE, and Train Pearson
# correlation in columns E:J are refreshed with the new run's numbers).

$ws.Range("E2").Value = -0.058
$ws.Range("F2").Value = -0.196
$ws.Range("G2").Value = 0.6860000000000001
$ws.Range("H2").Value = 0.828
$ws.Range("I2").Value = 0.62
$ws.Range("J2").Value = -0.013

$ws.Range("E3").Value = -0.074
$ws.Range("F3").Value = 1.69
$ws.Range("G3").Value = 7.41
$ws.Range("H3").Value = 2.722
$ws.Range("I3").Value = 2.26
$ws.Range("J3").Value = 0.042

$ws.Range("E4").Value = -0.006
$ws.Range("F4").Value = 1.647
$ws.Range("G4").Value = 0.639
$ws.Range("H4").Value = 0.799
$ws.Range("I4").Value = 0.57
$ws.Range("J4").Value = 0.748

$ws.Range("E5").Value = -0.027
$ws.Range("F5").Value = 1.66
$ws.Range("G5").Value = 0.179
$ws.Range("H5").Value = 0.423
$ws.Range("I5").Value = 0.32
$ws.Range("J5").Value = -0.117

$ws.Range("E6").Value = -0.479
$ws.Range("F6").Value = 1.951
$ws.Range("G6").Value = 1.162
$ws.Range("H6").Value = 1.078
$ws.Range("I6").Value = 0.842
$ws.Range("J6").Value = -0.282

$ws.Range("E7").Value = 0.005
$ws.Range("F7").Value = 1.64
$ws.Range("G7").Value = 0.882
$ws.Range("H7").Value = 0.9389999999999999
$ws.Range("I7").Value = 0.84
$ws.Range("J7").Value = 0.422

$ws.Range("E8").Value = -0.003
$ws.Range("F8").Value = 1.645
$ws.Range("G8").Value = 1.257
$ws.Range("H8").Value = 1.121
$ws.Range("I8").Value = 0.992
$ws.Range("J8").Value = 0.036

$ws.Range("E9").Value = -0.119
$ws.Range("F9").Value = 1.719
$ws.Range("G9").Value = 0.383
$ws.Range("H9").Value = 0.619
$ws.Range("I9").Value = 0.485
$ws.Range("J9").Value = -0.455

$ws.Range("E10").Value = -0.189
$ws.Range("F10").Value = 1.764
$ws.Range("G10").Value = 1.387
$ws.Range("H10").Value = 1.178
$ws.Range("I10").Value = 0.965
$ws.Range("J10").Value = 0.304

$ws.Range("E11").Value = -0.002
$ws.Range("F11").Value = 1.644
$ws.Range("G11").Value = 0.898
$ws.Range("H11").Value = 0.948
$ws.Range("I11").Value = 0.839
$ws.Range("J11").Value = 0.122

$ws.Range("E12").Value = 0.008999999999999999
$ws.Range("F12").Value = 1.637
$ws.Range("G12").Value = 0.625
$ws.Range("H12").Value = 0.791
$ws.Range("I12").Value = 0.678
$ws.Range("J12").Value = 0.442

$ws.Range("E13").Value = 0.008
$ws.Range("F13").Value = 1.638
$ws.Range("G13").Value = 0.3
$ws.Range("H13").Value = 0.548
$ws.Range("I13").Value = 0.432
$ws.Range("J13").Value = 0.345

$ws.Range("E14").Value = -0.07000000000000001
$ws.Range("F14").Value = 1.688
$ws.Range("G14").Value = 0.522
$ws.Range("H14").Value = 0.722
$ws.Range("I14").Value = 0.672
$ws.Range("J14").Value = 0.274

$ws.Range("E15").Value = -0.075
$ws.Range("F15").Value = 1.691
$ws.Range("G15").Value = 0.537
$ws.Range("H15").Value = 0.733
$ws.Range("I15").Value = 0.533
$ws.Range("J15").Value = -0.002

$ws.Range("E16").Value = -0.372
$ws.Range("F16").Value = 1.289
$ws.Range("G16").Value = 0.52
$ws.Range("H16").Value = 0.721
$ws.Range("I16").Value = 0.548
$ws.Range("J16").Value = -0.547

$ws.Range("E17").Value = -0.03
$ws.Range("F17").Value = 1.286
$ws.Range("G17").Value = 1.442
$ws.Range("H17").Value = 1.201
$ws.Range("I17").Value = 1.01
$ws.Range("J17").Value = 0.548

$ws.Range("E18").Value = -1.85
$ws.Range("F18").Value = 1.13
$ws.Range("G18").Value = 0.217
$ws.Range("H18").Value = 0.466
$ws.Range("I18").Value = 0.4
$ws.Range("J18").Value = 1

$ws.Range("E19").Value = -7.675
$ws.Range("F19").Value = 1.394
$ws.Range("G19").Value = 0.174
$ws.Range("H19").Value = 0.417
$ws.Range("I19").Value = 0.391
$ws.Range("J19").Value = -1

$ws.Range("E20").Value = -9.409000000000001
$ws.Range("F20").Value = 1.473
$ws.Range("G20").Value = 0.526
$ws.Range("H20").Value = 0.725
$ws.Range("I20").Value = 0.696
$ws.Range("J20").Value = 1

$ws.Range("E21").Value = -8.539
$ws.Range("F21").Value = 1.434
$ws.Range("G21").Value = 0.298
$ws.Range("H21").Value = 0.546
$ws.Range("I21").Value = 0.514
$ws.Range("J21").Value = -1

$ws.Range("E22").Value = -0.912
$ws.Range("F22").Value = 1.087
$ws.Range("G22").Value = 0.052
$ws.Range("H22").Value = 0.228
$ws.Range("I22").Value = 0.175
$ws.Range("J22").Value = 1

$ws.Range("E23").Value = -7.679
$ws.Range("F23").Value = 1.394
$ws.Range("G23").Value = 0.295
$ws.Range("H23").Value = 0.543
$ws.Range("I23").Value = 0.509
$ws.Range("J23").Value = -1

$ws.Range("E24").Value = -3.546
$ws.Range("F24").Value = 1.207
$ws.Range("G24").Value = 0.718
$ws.Range("H24").Value = 0.847
$ws.Range("I24").Value = 0.756
$ws.Range("J24").Value = 1

$ws.Range("E25").Value = -0.051
$ws.Range("F25").Value = 1.048
$ws.Range("G25").Value = 0.6889999999999999
$ws.Range("H25").Value = 0.83
$ws.Range("I25").Value = 0.8100000000000001
$ws.Range("J25").Value = -1

$ws.Range("E26").Value = -6.998
$ws.Range("F26").Value = 1.364
$ws.Range("G26").Value = 0.5679999999999999
$ws.Range("H26").Value = 0.754
$ws.Range("I26").Value = 0.6929999999999999
$ws.Range("J26").Value = -1

$ws.Range("E27").Value = -0.226
$ws.Range("F27").Value = 1.056
$ws.Range("G27").Value = 0.013
$ws.Range("H27").Value = 0.114
$ws.Range("I27").Value = 0.148
$ws.Range("J27").Value = 1

$ws.Range("E28").Value = -0.033
$ws.Range("F28").Value = 1.047
$ws.Range("G28").Value = 1.307
$ws.Range("H28").Value = 1.143
$ws.Range("I28").Value = 1.125
$ws.Range("J28").Value = -1

$ws.Range("E29").Value = -2.586
$ws.Range("F29").Value = 1.163
$ws.Range("G29").Value = 0.553
$ws.Range("H29").Value = 0.744
$ws.Range("I29").Value = 0.639
$ws.Range("J29").Value = 1

$ws.Range("E30").Value = -1461.147
$ws.Range("F30").Value = 67.461
$ws.Range("G30").Value = 1.711
$ws.Range("H30").Value = 1.308
$ws.Range("I30").Value = 1.308
$ws.Range("J30").Value = 1

$ws.Range("E31").Value = -0.036
$ws.Range("F31").Value = 1.047
$ws.Range("G31").Value = 0.079
$ws.Range("H31").Value = 0.281
$ws.Range("I31").Value = 0.276
$ws.Range("J31").Value = 1

$ws.Range("E32").Value = -24.494
$ws.Range("F32").Value = 2.159
$ws.Range("G32").Value = 2.267
$ws.Range("H32").Value = 1.506
$ws.Range("I32").Value = 1.475
$ws.Range("J32").Value = -1

# --- Summary rows 33-36: TotalNNAvg now leads the block, followed by the
# Pre2020 / Trans / Post2020 averages (each shifted down one row), all with
# refreshed average values ---

$ws.Range("A33").Value = "TotalNNAvg"
$ws.Range("E33").Value = -11.44310609243697
$ws.Range("F33").Value = 1.933707983193277
$ws.Range("G33").Value = 0.9418518907563024
$ws.Range("H33").Value = 0.8628266806722691
$ws.Range("I33").Value = 0.7057090336134454
$ws.Range("J33").Value = -0.07109663865546222

$ws.Range("A34").Value = "Pre2020NNavg"
$ws.Range("E34").Value = -0.07714285714285714
$ws.Range("F34").Value = 1.558428571428571
$ws.Range("G34").Value = 1.204785714285714
$ws.Range("H34").Value = 0.9606428571428572
$ws.Range("I34").Value = 0.7891428571428571
$ws.Range("J34").Value = 0.1332857142857143

$ws.Range("A35").Value = "TransNNavg"
$ws.Range("E35").Value = -0.372
$ws.Range("F35").Value = 1.289
$ws.Range("G35").Value = 0.52
$ws.Range("H35").Value = 0.721
$ws.Range("I35").Value = 0.548
$ws.Range("J35").Value = -0.547

$ws.Range("A36").Value = "Post2020NNavg"
$ws.Range("E36").Value = -95.9506875
$ws.Range("F36").Value = 5.421875
$ws.Range("G36").Value = 0.6818124999999999
$ws.Range("H36").Value = 0.7283125
$ws.Range("I36").Value = 0.6828124999999999
$ws.Range("J36").Value = 0.09675
